# Attendance2022.xlsx update:
#  - The tracked user's e-mail address changed from
#    anita24feb.sinha@gmail.com -> avinash18dce.kumar@gmail.com
#    on every month sheet (cell A3).
#  - On the "May" sheet, day 27 (column AB) attendance was marked
#    "P" (present) for both tracked rows (AB2 and AB3).

$wb = $excel.ActiveWorkbook

$monthNames = @("January","February","March","April","May","June","July","August","September","October","November","December")

foreach ($name in $monthNames) {
    $ws = $wb.Worksheets.Item($name)

    # Update the e-mail address in A3 on every month sheet.
    $ws.Range("A3").Value = "avinash18dce.kumar@gmail.com"
}

# May sheet gets the extra attendance marks for day 27 (column AB).
$may = $wb.Worksheets.Item("May")
$may.Range("AB2").Value = "P"
$may.Range("AB3").Value = "P"
